# Reorder slides: move the "Indian Socio-Cultural Context" slide
# (currently slide 18) so it lands after the two picture-only slides
# that follow it, i.e. to position 20.
#
# Before: ... 17, [18: Indian Socio-Cultural Context], [19: picture], [20: picture], 21, ...
# After:  ... 17, [18: picture], [19: picture], [20: Indian Socio-Cultural Context], 21, ...

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$s.MoveTo(20)
